$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Row 2 updates
$ws.Range("K2").Value = 3
$ws.Range("L2").Value = 1
$ws.Range("M2").Value = 0.165415
$ws.Range("N2").Value = 0.496245
$ws.Range("O2").Value = 0.06044768156291203
$ws.Range("P2").Value = 0.06044768156291203
$ws.Range("Q2").Value = 0.028304656895
$ws.Range("R2").Value = 0.254741912055
$ws.Range("S2").Value = 0.06044768156291203
$ws.Range("T2").Value = 0.06044768156291203

# Row 3 updates
$ws.Range("O3").Value = 0.8243344049378915
$ws.Range("P3").Value = 0.8243344049378915
$ws.Range("S3").Value = 0.8243344049378915
$ws.Range("T3").Value = 0.8243344049378915

# Row 4 updates
$ws.Range("O4").Value = 0.1152179134991965
$ws.Range("P4").Value = 0.1152179134991965
$ws.Range("S4").Value = 0.1152179134991965
$ws.Range("T4").Value = 0.1152179134991965
